$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04155
$ws.Range("H2").Value = 0.12465
$ws.Range("I2").Value = 0.0001466168179836329
$ws.Range("J2").Value = 0.0001466168179836329
$ws.Range("M2").Value = 5.008808666666667
$ws.Range("N2").Value = 15.026426
$ws.Range("O2").Value = 0.3739112966508367
$ws.Range("P2").Value = 0.3739112966508367
$ws.Range("Q2").Value = 0.2081160001
$ws.Range("R2").Value = 1.8730440009
$ws.Range("S2").Value = [double]"5.482168452307987E-05"
$ws.Range("T2").Value = [double]"5.482168452307987E-05"

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04155
$ws.Range("H3").Value = 0.12465
$ws.Range("I3").Value = 0.0001466168179836329
$ws.Range("J3").Value = 0.0001466168179836329
$ws.Range("O3").Value = 0.3132668979860996
$ws.Range("P3").Value = 0.3132668979860996
$ws.Range("Q3").Value = 0.17436181885
$ws.Range("R3").Value = 1.56925636965
$ws.Range("S3").Value = [double]"4.593019576232524E-05"
$ws.Range("T3").Value = [double]"4.593019576232526E-05"

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.04155
$ws.Range("H4").Value = 0.12465
$ws.Range("I4").Value = 0.0001466168179836329
$ws.Range("J4").Value = 0.0001466168179836329
$ws.Range("M4").Value = 4.190471333333333
$ws.Range("O4").Value = 0.3128218053630638
$ws.Range("P4").Value = 0.3128218053630638
$ws.Range("Q4").Value = 0.1741140839
$ws.Range("R4").Value = 1.5670267551
$ws.Range("S4").Value = [double]"4.586493769822775E-05"
$ws.Range("T4").Value = [double]"4.586493769822775E-05"

# Row 5
$ws.Range("I5").Value = 0.9992428949822291
$ws.Range("J5").Value = 0.9992428949822291
$ws.Range("M5").Value = 5.008808666666667
$ws.Range("N5").Value = 15.026426
$ws.Range("O5").Value = 0.3739112966508367
$ws.Range("P5").Value = 0.3739112966508367
$ws.Range("Q5").Value = 1418.380492033736
$ws.Range("R5").Value = 12765.42442830362
$ws.Range("S5").Value = 0.3736282065319411
$ws.Range("T5").Value = 0.3736282065319411

# Row 6
$ws.Range("I6").Value = 0.9992428949822291
$ws.Range("J6").Value = 0.9992428949822291
$ws.Range("O6").Value = 0.3132668979860996
$ws.Range("P6").Value = 0.3132668979860996
$ws.Range("S6").Value = 0.3130297220457328
$ws.Range("T6").Value = 0.3130297220457328

# Row 7
$ws.Range("I7").Value = 0.9992428949822291
$ws.Range("J7").Value = 0.9992428949822291
$ws.Range("M7").Value = 4.190471333333333
$ws.Range("O7").Value = 0.3128218053630638
$ws.Range("P7").Value = 0.3128218053630638
$ws.Range("Q7").Value = 1186.646004504317
$ws.Range("S7").Value = 0.3125849664045552
$ws.Range("T7").Value = 0.3125849664045552

# Row 8
$ws.Range("I8").Value = 0.0006104881997874136
$ws.Range("J8").Value = 0.0006104881997874135
$ws.Range("M8").Value = 5.008808666666667
$ws.Range("N8").Value = 15.026426
$ws.Range("O8").Value = 0.3739112966508367
$ws.Range("P8").Value = 0.3739112966508367
$ws.Range("Q8").Value = 0.8665606305968888
$ws.Range("R8").Value = 7.799045675372001
$ws.Range("S8").Value = 0.0002282684343725468
$ws.Range("T8").Value = 0.0002282684343725468

# Row 9
$ws.Range("I9").Value = 0.0006104881997874136
$ws.Range("J9").Value = 0.0006104881997874135
$ws.Range("O9").Value = 0.3132668979860996
$ws.Range("P9").Value = 0.3132668979860996
$ws.Range("R9").Value = 6.534124183622001
$ws.Range("S9").Value = 0.0001912457446045213
$ws.Range("T9").Value = 0.0001912457446045213

# Row 10
$ws.Range("I10").Value = 0.0006104881997874136
$ws.Range("J10").Value = 0.0006104881997874135
$ws.Range("M10").Value = 4.190471333333333
$ws.Range("O10").Value = 0.3128218053630638
$ws.Range("P10").Value = 0.3128218053630638
$ws.Range("Q10").Value = 0.7249822707897776
$ws.Range("R10").Value = 6.524840437107999
$ws.Range("S10").Value = 0.0001909740208103455
$ws.Range("T10").Value = 0.0001909740208103454

